$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # MatchCase = $true, MatchWholeWord = $true so that short generic
    # words (e.g. "Riddle", "Solution", "Date") only replace the exact
    # table-label runs and don't clobber case-insensitive / partial-word
    # occurrences buried inside unrelated sentences elsewhere in the doc.
    $result = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                                       $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "WARNING: replace failed for '$old'"
    }
}

# Longer / more specific strings must be replaced before any shorter
# strings that are substrings of them (e.g. "General VMC Video
# Introduction" contains "Video Introduction").
Replace-Text "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"

Replace-Text "Video Title" "Kichwa cha Video"
Replace-Text "Topic" "Mada"
Replace-Text "Aim(s)" "Malengo"
Replace-Text "Length" "Urefu"
Replace-Text "Camp Location" "Mahali pa Kambi"
Replace-Text "Facilitators" "Wawezeshaji"
Replace-Text "N. of students" "N. ya wanafunzi"
Replace-Text "Date" "Tarehe"
Replace-Text "Resources" "Rasilimali"
Replace-Text "needed" "inahitajika"
Replace-Text "Preparations" "Maandalizi"
Replace-Text "Video time" "Muda wa video"
Replace-Text "What facilitator does" "Mwezeshaji anafanya nini"
Replace-Text "What learners do" "Wanachofanya wanafunzi"
Replace-Text "Video Introduction" "Utangulizi wa Video"
Replace-Text "Riddle" "Kitendawili"
Replace-Text "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"
Replace-Text "Solution" "Suluhisho"

# Update the document's default language (Swahili: Tanzania -> Kenya).
$d.Styles("Normal").LanguageID = "sw-KE"

Write-Output "done"
